$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New version entry (0.2.0) is added as row 12, below the existing last
# row (11). Start by copying the formatting of the previous "new row"
# (row 10) down onto row 12, which reproduces the alternating-row fill
# / border / alignment / number-format styling used throughout the
# table (this matches s="4" for columns A-D and s="8" for columns E-G).
$ws.Range("A10:G10").Copy($ws.Range("A12:G12"))

# Match the row height used by the other wrapped-text rows of this kind.
$ws.Rows.Item(12).RowHeight = $ws.Rows.Item(10).RowHeight

# Fill in the actual values for the new release (the change log text is
# entered before the open points text so the two new strings land in the
# shared-string table in the same order Excel produced them in).
$ws.Range("A12").Value = "0.2.0"
$ws.Range("B12").Value = "AUTOMATA CELULAR - copia (16)"
$ws.Range("D12").Value = "-Changed reproduction and distribution to two parts.`n-Disaggregate done in between reproduction and distribution.`n-Aggregation conditions.`n-Fixed minor error on the SG.`n-With 4 or less niches the distribution is now equaly done."
$ws.Range("C12").Value = "-Make the code more readable using functions.`n-UI: Delete rows according to working functionality.`n*Implement mutations.`n*Graphic representation of F'.`n-Rework E calc.`n-Document every function."
$ws.Range("E12").Value = "Python 3.6.1"
$ws.Range("F12").Value = "Qt version: 5.6.2`nSIP version: 4.18`nPyQt version: 5.6"
$ws.Range("G12").Value = " PyInstaller 3.3.1"

# Move the active selection down to the newly added row, the same way
# Excel records where the user ended up after adding the row.
$ws.Activate()
$ws.Range("C13").Select()
